$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "training" and add a new "test" sheet right
# after it. Worksheets.Add() inserts *before* the sheet passed as "After"
# when not given explicitly, so pass the "training" sheet as After to land
# "test" in the second slot (matching sheetId=1/training, sheetId=2/test).
$wb.Worksheets.Item(1).Name = "training"
$wsNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item(1))
$wsNew.Name = "test"

function Fill-Stats {
    param($ws, $counts)

    $ws.Range("A1").Value = "Class"
    $ws.Range("B1").Value = "Count"
    $ws.Range("C1").Value = "Percentage"

    for ($i = 0; $i -lt 5; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i
        $ws.Cells.Item($row, 2).Value = $counts[$i]
        $ws.Cells.Item($row, 3).Formula = "=B$row/SUM(B2:B6)"
        $ws.Cells.Item($row, 3).NumberFormat = "0.0%"
    }
}

Fill-Stats $wb.Worksheets.Item(1) @(25810, 2443, 5292, 873, 708)
Fill-Stats $wb.Worksheets.Item(2) @(39533, 3762, 7861, 1214, 1206)

# training: widen the Percentage column (and the following, empty column)
# the way Excel's "AutoFit selection" leaves things after the user eyeballs it.
$wb.Worksheets.Item(1).Columns.Item(3).ColumnWidth = 10.5546875
$wb.Worksheets.Item(1).Columns.Item(4).ColumnWidth = 11.44140625

# Restore the selections/active sheet the workbook was left on.
$wb.Worksheets.Item(2).Range("B6").Select() | Out-Null
$wb.Worksheets.Item(1).Select() | Out-Null
$wb.Worksheets.Item(1).Range("D15").Select() | Out-Null
